$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Distance" column (L) header
$ws.Range("L1").Value = "Distance"

# Per-row distance values (row 2 .. row 57)
$distances = @(
    101,
    22.6,
    146,
    83,
    124,
    11.5,
    11.5,
    207,
    68.9,
    140,
    140,
    97.8,
    130,
    55.9,
    55.9,
    112,
    112,
    143,
    283,
    201,
    201,
    38.5,
    65.9,
    27.5,
    1.8,
    47.1,
    1.3,
    81.8,
    87.7,
    87.7,
    39.7,
    39.7,
    107,
    21.8,
    80.7,
    80.7,
    1.4,
    1.2,
    110,
    110,
    1.3,
    156,
    156,
    115,
    123,
    123,
    287,
    287,
    179,
    179,
    150,
    150,
    67.8,
    186,
    132,
    132
)

for ($i = 0; $i -lt $distances.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $distances[$i]
}

# Reflect the post-edit view state: scrolled down, with the next empty
# cell below the new column selected.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("L58").Select()
